$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.833.75'
$ws.Range("E2").Value = '  -0.78%  '

$ws.Range("D3").Value = '3.519.08'
$ws.Range("E3").Value = '  -1.47%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.38'
$ws.Range("E5").Value = '  +2.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.83'
$ws.Range("E6").Value = '  -1.86%  '

$ws.Range("D7").Value = '3.513.22'
$ws.Range("E7").Value = '  -1.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.611'
$ws.Range("E8").Value = '  -1.72%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.210'
$ws.Range("E10").Value = '  +4.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.637'
$ws.Range("E11").Value = '  -1.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.19'
$ws.Range("E12").Value = '  -2.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000304'
$ws.Range("E13").Value = '  -1.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.37'
$ws.Range("E14").Value = '  -2.35%  '

$ws.Range("D15").Value = '4.091.59'
$ws.Range("E15").Value = '  -1.14%  '

$ws.Range("D16").Value = '69.912.28'
$ws.Range("E16").Value = '  -0.68%  '

$ws.Range("D17").Value = '3.576.50'
$ws.Range("E17").Value = '  +0.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.76'
$ws.Range("E18").Value = '  -3.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.52'
$ws.Range("E19").Value = '  +0.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '573.50'
$ws.Range("E20").Value = '  +3.71%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.982'
$ws.Range("E22").Value = '  -3.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.18'
$ws.Range("E23").Value = '  -4.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.65'
$ws.Range("E24").Value = '  -0.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.79'
$ws.Range("E25").Value = '  -2.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '93.46'
$ws.Range("E26").Value = '  -2.95%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("E27").Value = '  -5.09%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.91'
$ws.Range("E28").Value = '  -2.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  +0.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.82'
$ws.Range("E30").Value = '  -1.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.96'
$ws.Range("E31").Value = '  -5.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.11'
$ws.Range("E32").Value = '  -3.73%  '

$ws.Range("E33").Value = '  -1.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '62.95'
$ws.Range("E34").Value = '  -3.42%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.29'
$ws.Range("E35").Value = '  +1.21%  '

$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.59'
$ws.Range("E36").Value = '  +14.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '530.13'
$ws.Range("E37").Value = '  -4.85%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.399'
$ws.Range("E39").Value = '  -4.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.70'
$ws.Range("E40").Value = '  -4.23%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.525.90'
$ws.Range("E41").Value = '  +4.28%  '

$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0772'
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.48'
$ws.Range("E43").Value = '  +2.67%  '

$ws.Range("E44").Value = '  -0.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0452'
$ws.Range("E45").Value = '  +0.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.41'
$ws.Range("E46").Value = '  -5.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.88'
$ws.Range("E47").Value = '  -3.87%  '

$ws.Range("E48").Value = '  +1.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.07'
$ws.Range("E49").Value = '  -1.50%  '

$ws.Range("E50").Value = '  +0.38%  '

$ws.Range("E51").Value = '  -4.56%  '
